# Applies the natmiOut Tnc-Itgav LR-pair edit: "Natmi following Dr Hou advice"
# Re-computed rows for Sending cluster x Target cluster pairs across ECs/FAPs/sCs
# (Ligand=Tnc, Receptor=Itgav), extending the table from 6 to 9 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnc"
$ws.Cells.Item(2, 3).Value = "Itgav"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.442371333333333
$ws.Cells.Item(2, 8).Value = 4.327114
$ws.Cells.Item(2, 9).Value = 0.03522044016446201
$ws.Cells.Item(2, 10).Value = 0.03522044016446201
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 16.535604
$ws.Cells.Item(2, 14).Value = 49.606812
$ws.Cells.Item(2, 15).Value = 0.2120453146491552
$ws.Cells.Item(2, 16).Value = 0.2120453146491552
$ws.Cells.Item(2, 17).Value = 23.850481188952
$ws.Cells.Item(2, 18).Value = 214.654330700568
$ws.Cells.Item(2, 19).Value = 0.007468329316755091
$ws.Cells.Item(2, 20).Value = 0.007468329316755093

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnc"
$ws.Cells.Item(3, 3).Value = "Itgav"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.442371333333333
$ws.Cells.Item(3, 8).Value = 4.327114
$ws.Cells.Item(3, 9).Value = 0.03522044016446201
$ws.Cells.Item(3, 10).Value = 0.03522044016446201
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 40.62063066666667
$ws.Cells.Item(3, 14).Value = 121.861892
$ws.Cells.Item(3, 15).Value = 0.5209011059384622
$ws.Cells.Item(3, 16).Value = 0.5209011059384622
$ws.Cells.Item(3, 17).Value = 58.59003321552089
$ws.Cells.Item(3, 18).Value = 527.310298939688
$ws.Cells.Item(3, 19).Value = 0.0183463662333077
$ws.Cells.Item(3, 20).Value = 0.0183463662333077

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnc"
$ws.Cells.Item(4, 3).Value = "Itgav"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.442371333333333
$ws.Cells.Item(4, 8).Value = 4.327114
$ws.Cells.Item(4, 9).Value = 0.03522044016446201
$ws.Cells.Item(4, 10).Value = 0.03522044016446201
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 20.825229
$ws.Cells.Item(4, 14).Value = 62.475687
$ws.Cells.Item(4, 15).Value = 0.2670535794123827
$ws.Cells.Item(4, 16).Value = 0.2670535794123827
$ws.Cells.Item(4, 17).Value = 30.037713319702
$ws.Cells.Item(4, 18).Value = 270.339419877318
$ws.Cells.Item(4, 19).Value = 0.009405744614399228
$ws.Cells.Item(4, 20).Value = 0.009405744614399228

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnc"
$ws.Cells.Item(5, 3).Value = "Itgav"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.331397
$ws.Cells.Item(5, 8).Value = 33.994191
$ws.Cells.Item(5, 9).Value = 0.2766948987373093
$ws.Cells.Item(5, 10).Value = 0.2766948987373092
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 16.535604
$ws.Cells.Item(5, 14).Value = 49.606812
$ws.Cells.Item(5, 15).Value = 0.2120453146491552
$ws.Cells.Item(5, 16).Value = 0.2120453146491552
$ws.Cells.Item(5, 17).Value = 187.371493558788
$ws.Cells.Item(5, 18).Value = 1686.343442029092
$ws.Cells.Item(5, 19).Value = 0.05867185686456889
$ws.Cells.Item(5, 20).Value = 0.05867185686456888

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnc"
$ws.Cells.Item(6, 3).Value = "Itgav"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.331397
$ws.Cells.Item(6, 8).Value = 33.994191
$ws.Cells.Item(6, 9).Value = 0.2766948987373093
$ws.Cells.Item(6, 10).Value = 0.2766948987373092
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 40.62063066666667
$ws.Cells.Item(6, 14).Value = 121.861892
$ws.Cells.Item(6, 15).Value = 0.5209011059384622
$ws.Cells.Item(6, 16).Value = 0.5209011059384622
$ws.Cells.Item(6, 17).Value = 460.2884924743747
$ws.Cells.Item(6, 18).Value = 4142.596432269373
$ws.Cells.Item(6, 19).Value = 0.1441306787597952
$ws.Cells.Item(6, 20).Value = 0.1441306787597952

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnc"
$ws.Cells.Item(7, 3).Value = "Itgav"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.331397
$ws.Cells.Item(7, 8).Value = 33.994191
$ws.Cells.Item(7, 9).Value = 0.2766948987373093
$ws.Cells.Item(7, 10).Value = 0.2766948987373092
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 20.825229
$ws.Cells.Item(7, 14).Value = 62.475687
$ws.Cells.Item(7, 15).Value = 0.2670535794123827
$ws.Cells.Item(7, 16).Value = 0.2670535794123827
$ws.Cells.Item(7, 17).Value = 235.978937414913
$ws.Cells.Item(7, 18).Value = 2123.810436734217
$ws.Cells.Item(7, 19).Value = 0.07389236311294521
$ws.Cells.Item(7, 20).Value = 0.0738923631129452

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Tnc"
$ws.Cells.Item(8, 3).Value = "Itgav"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 28.17890933333333
$ws.Cells.Item(8, 8).Value = 84.536728
$ws.Cells.Item(8, 9).Value = 0.6880846610982287
$ws.Cells.Item(8, 10).Value = 0.6880846610982286
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 16.535604
$ws.Cells.Item(8, 14).Value = 49.606812
$ws.Cells.Item(8, 15).Value = 0.2120453146491552
$ws.Cells.Item(8, 16).Value = 0.2120453146491552
$ws.Cells.Item(8, 17).Value = 465.955285887904
$ws.Cells.Item(8, 18).Value = 4193.597572991136
$ws.Cells.Item(8, 19).Value = 0.1459051284678312
$ws.Cells.Item(8, 20).Value = 0.1459051284678312

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Tnc"
$ws.Cells.Item(9, 3).Value = "Itgav"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 28.17890933333333
$ws.Cells.Item(9, 8).Value = 84.536728
$ws.Cells.Item(9, 9).Value = 0.6880846610982287
$ws.Cells.Item(9, 10).Value = 0.6880846610982286
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 40.62063066666667
$ws.Cells.Item(9, 14).Value = 121.861892
$ws.Cells.Item(9, 15).Value = 0.5209011059384622
$ws.Cells.Item(9, 16).Value = 0.5209011059384622
$ws.Cells.Item(9, 17).Value = 1144.64506861882
$ws.Cells.Item(9, 18).Value = 10301.80561756938
$ws.Cells.Item(9, 19).Value = 0.3584240609453593
$ws.Cells.Item(9, 20).Value = 0.3584240609453592

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tnc"
$ws.Cells.Item(10, 3).Value = "Itgav"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 28.17890933333333
$ws.Cells.Item(10, 8).Value = 84.536728
$ws.Cells.Item(10, 9).Value = 0.6880846610982287
$ws.Cells.Item(10, 10).Value = 0.6880846610982286
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 20.825229
$ws.Cells.Item(10, 14).Value = 62.475687
$ws.Cells.Item(10, 15).Value = 0.2670535794123827
$ws.Cells.Item(10, 16).Value = 0.2670535794123827
$ws.Cells.Item(10, 17).Value = 586.832239836904
$ws.Cells.Item(10, 18).Value = 5281.490158532136
$ws.Cells.Item(10, 19).Value = 0.1837554716850383
$ws.Cells.Item(10, 20).Value = 0.1837554716850382
